# Refresh the "Enterprises density", "Employment (% of total)",
# "Enterprises (% of total)" and "Value added to the economy (% of total)"
# statistics on the Poland Summary sheet with more precise (two-decimal)
# figures, per the 2015-04-01 autogenerated data refresh.
#
# These figures are stored as text (not numbers) in the sheet, so each
# target cell is explicitly formatted as Text before the new value is
# written, preserving the original "text-number" representation instead
# of letting Excel auto-coerce the digits into a numeric cell. The number
# format is restored afterwards so the cells keep displaying/behaving like
# the rest of the (General-formatted) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B13" = "37.68"   # Enterprises density (per 1000 people) - Micro
    "C13" = "1.77"    # Enterprises density (per 1000 people) - SMEs
    "D13" = "39.46"   # Enterprises density (per 1000 people) - MSMEs

    "B14" = "36.19"   # Employment (% of total) - Micro
    "C14" = "32.54"   # Employment (% of total) - SMEs
    "D14" = "68.73"   # Employment (% of total) - MSMEs

    "B16" = "95.31"   # Enterprises (% of total) - Micro
    "C16" = "4.49"    # Enterprises (% of total) - SMEs

    "B20" = "16.42"   # Value added to the economy (% of total) - Micro
    "C20" = "34.88"   # Value added to the economy (% of total) - SMEs
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.NumberFormat = "General"
}
